$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.783.98'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = '1.616.54'
$ws.Range('E3').Value = '  -1.50%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.38'
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.520'
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.19'
$ws.Range('E8').Value = '  -1.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.256'
$ws.Range('E9').Value = '  -1.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0608'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0877'
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('D12').Value = '1.840.96'
$ws.Range('E12').Value = '  -1.74%  '
$ws.Range('D13').Value = '1.617.67'
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.00'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.558'
$ws.Range('E15').Value = '  -3.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.92'
$ws.Range('E16').Value = '  -1.51%  '
$ws.Range('D17').Value = '27.738.80'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.33'
$ws.Range('E18').Value = '  -2.77%  '
$ws.Range('D19').Value = '0.0₃0718'
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('E20').Value = '  -0.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.996'
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.28'
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.07'
$ws.Range('E23').Value = '  -5.88%  '
$ws.Range('E24').Value = '  -3.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.21'
$ws.Range('E25').Value = '  +1.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.86'
$ws.Range('E26').Value = '  -1.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.111'
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.44'
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0478'
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('E33').Value = '  -2.01%  '
$ws.Range('D34').Value = '1.385.71'
$ws.Range('E34').Value = '  -2.56%  '
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.984'
$ws.Range('E36').Value = '  +8.74%  '
$ws.Range('E37').Value = '  -1.33%  '
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('E39').Value = '  -1.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.849'
$ws.Range('E40').Value = '  -3.88%  '
$ws.Range('E41').Value = '  -1.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.996'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.82'
$ws.Range('E43').Value = '  -2.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.27'
$ws.Range('E44').Value = '  -2.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.43'
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.15'
$ws.Range('E46').Value = '  -2.52%  '
$ws.Range('D47').Value = '1.750.90'
$ws.Range('E47').Value = '  -1.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.47'
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0502'
$ws.Range('E50').Value = '  -0.88%  '
$ws.Range('D51').Value = '0.0₇0967'
$ws.Range('E51').Value = '  -8.27%  '
